$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": August sales for PORCELANATO (col M)
# newly recorded for three clients.
# ---------------------------------------------------------------
$wsGrupo.Range("M16").Value = -23.16
$wsGrupo.Range("M22").Value = 1634.69
$wsGrupo.Range("M36").Value = 8963.41

# Footer counter row: number of non-zero PORCELANATO entries out of 54
$wsGrupo.Range("M56").Value = "3 de 54"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL": same amounts land in the "agosto" column (F)
# for the same three clients, plus the column total.
# ---------------------------------------------------------------
$wsMensual.Range("F16").Value = -23.16
$wsMensual.Range("F22").Value = 1634.69
$wsMensual.Range("F36").Value = 8963.41
$wsMensual.Range("F56").Value = 10448.22

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": VENTA (D) per product group is the sum
# of that group's column on "VENTAS POR GRUPO"; POR CUMPLIR (E) and
# CUMPLIMIENTO (F) are recomputed, plus a few direct PRESUPUESTO (C)
# budget updates.
# ---------------------------------------------------------------

# 240X120 PORCELANATO
$wsCumpl.Range("D2").Value = 0
$wsCumpl.Range("E2").Value = 9970.34304517915
$wsCumpl.Range("F2").Value = 0

# 240X80 PORCELANATO
$wsCumpl.Range("C3").Value = 27457.0076
$wsCumpl.Range("D3").Value = 380.16
$wsCumpl.Range("E3").Value = 27076.8476
$wsCumpl.Range("F3").Value = 0.01384564572870643

# FREGADEROS DE COCINA
$wsCumpl.Range("D4").Value = 0
$wsCumpl.Range("E4").Value = 1003
$wsCumpl.Range("F4").Value = 0

# GRANITO
$wsCumpl.Range("D5").Value = 0
$wsCumpl.Range("E5").Value = 238.32
$wsCumpl.Range("F5").Value = 0

# INODOROS
$wsCumpl.Range("D7").Value = 0
$wsCumpl.Range("E7").Value = 2400
$wsCumpl.Range("F7").Value = 0

# LAVABOS
$wsCumpl.Range("D8").Value = 0
$wsCumpl.Range("E8").Value = 1000
$wsCumpl.Range("F8").Value = 0

# LED
$wsCumpl.Range("D9").Value = 0
$wsCumpl.Range("E9").Value = 300
$wsCumpl.Range("F9").Value = 0

# NO RESURTIBLES
$wsCumpl.Range("D10").Value = 0
$wsCumpl.Range("E10").Value = 1300.5
$wsCumpl.Range("F10").Value = 0

# PANELES DECORATIVOS (budget revised 350 -> 100)
$wsCumpl.Range("C12").Value = 100
$wsCumpl.Range("D12").Value = 0
$wsCumpl.Range("E12").Value = 100
$wsCumpl.Range("F12").Value = 0

# PANELES PU (budget revised 130 -> 20)
$wsCumpl.Range("C13").Value = 20
$wsCumpl.Range("E13").Value = 20

# PANELES PVC (budget revised 240 -> 100)
$wsCumpl.Range("C14").Value = 100
$wsCumpl.Range("E14").Value = 100

# PIEDRA SINTERIZADA
$wsCumpl.Range("D15").Value = -644.89
$wsCumpl.Range("E15").Value = 14144.89
$wsCumpl.Range("F15").Value = -0.04776962962962963

# PORCELANATO (budget revised 51826.46 -> 56059.7)
$wsCumpl.Range("C16").Value = 56059.7
$wsCumpl.Range("D16").Value = 10712.95
$wsCumpl.Range("E16").Value = 45346.75
$wsCumpl.Range("F16").Value = 0.1910989534371394

# PUERTAS DE SEGURIDAD
$wsCumpl.Range("D17").Value = 0
$wsCumpl.Range("E17").Value = 684
$wsCumpl.Range("F17").Value = 0

# TOTAL row
$wsCumpl.Range("C19").Value = 117439.6906451791
$wsCumpl.Range("D19").Value = 10448.22
$wsCumpl.Range("E19").Value = 106991.4706451792
$wsCumpl.Range("F19").Value = 0.08896668530545808

# Column widths on "CUMPLIMIENTO MENSUAL" narrow slightly for "POR CUMPLIR"
# (col E) and widen for "CUMPLIMIENTO" (col F) once the new, shorter /
# longer formatted numbers are in place. The engine's ColumnWidth setter
# adds a fixed ~0.8333 character padding before it round-trips through
# xlsx, so back that out to land on the exact target widths.
$wsCumpl.Columns.Item(5).ColumnWidth = 21.166666666666668
$wsCumpl.Columns.Item(6).ColumnWidth = 25.166666666666668
